$wb = $excel.ActiveWorkbook

# --- Sheet "Ranges": update Sarnia's (row 9) Min values and Max values,
#     and Trenton's (row 3) Max values, to reflect the new final Wind data. ---
$wsRanges = $wb.Worksheets.Item("Ranges")
$wsRanges.Range("D3").Value = 21541
$wsRanges.Range("E3").Value = 154

$wsRanges.Range("A9").Value = 38928
$wsRanges.Range("B9").Value = 65
$wsRanges.Range("D9").Value = 33324
$wsRanges.Range("E9").Value = 159

# --- Sheet "Minimum": table is sorted ascending by Min_SpdOfMaxGust (col B).
#     Sarnia's min speed dropped from 78 to 65, so it now sorts in right
#     after Trenton (65) and before Toronto Airport (69); the rows that used
#     to hold Toronto Airport / Wiarton / Hamilton shift down by one, and the
#     old last Sarnia row becomes Hamilton. ---
$wsMin = $wb.Worksheets.Item("Minimum")

$wsMin.Range("A5").Value = 38928
$wsMin.Range("B5").Value = 65
$wsMin.Range("C5").Value = "Sarnia"

$wsMin.Range("A6").Value = 43123
$wsMin.Range("B6").Value = 69
$wsMin.Range("C6").Value = "Toronto Airport"

$wsMin.Range("A7").Value = 43126
$wsMin.Range("B7").Value = 72
$wsMin.Range("C7").Value = "Wiarton"

$wsMin.Range("C8").Value = "Hamilton"

# --- Sheet "Maximum": table is sorted descending by Max_SpdOfMaxGust (col B).
#     Sarnia's max speed rose to 159 and Trenton's rose to 154, so Sarnia now
#     sorts first and Trenton second; London / Hamilton shift down by one. ---
$wsMax = $wb.Worksheets.Item("Maximum")

$wsMax.Range("A2").Value = 33324
$wsMax.Range("B2").Value = 159
$wsMax.Range("C2").Value = "Sarnia"

$wsMax.Range("A3").Value = 21541
$wsMax.Range("B3").Value = 154

$wsMax.Range("A4").Value = 33772
$wsMax.Range("B4").Value = 148
$wsMax.Range("C4").Value = "London"

$wsMax.Range("C5").Value = "Hamilton"
